$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Workbook / sheet structure -------------------------------------------
# Original sheets: Tabelle1 (colorscheme "ovgu" data), Tabelle2 (empty), Tabelle3 (empty)
$shExtra = $wb.Worksheets.Item("Tabelle3")

# Drop the unused third sheet
$shExtra.Delete()

# Put the (currently empty) sheet that will hold the new "uulm" data first.
# NOTE: after .Move() the original object handles can refer to the swapped
# positions, so re-fetch the sheets by name afterwards to be safe.
$wb.Worksheets.Item("Tabelle2").Move($wb.Worksheets.Item("Tabelle1"))

$shUulm = $wb.Worksheets.Item("Tabelle2")
$shOvgu = $wb.Worksheets.Item("Tabelle1")

# Rename sheets to their final names
$shUulm.Name = "uulm"
$shOvgu.Name = "ovgu"

# --- Build the "uulm" colorscheme sheet ------------------------------------
$ws = $shUulm

# Row 1: parameters (C1 = base percentage 100%, G1 = mix-with-white percentage 30%)
$ws.Range("C1:F1").Merge()
$ws.Range("C1").Value = 1
$ws.Range("C1:F1").NumberFormat = "0%"
$ws.Range("C1:F1").HorizontalAlignment = -4108

$ws.Range("G1:J1").Merge()
$ws.Range("G1").Value = 0.3
$ws.Range("G1:J1").NumberFormat = "0%"
$ws.Range("G1:J1").HorizontalAlignment = -4108

# Row 2: column headers
$ws.Range("A2").NumberFormat = "0%"

$ws.Range("C2").Value = "r"
$ws.Range("D2").Value = "g"
$ws.Range("E2").Value = "b"
$ws.Range("F2").Value = "html"
$ws.Range("G2").Value = "r"
$ws.Range("H2").Value = "g"
$ws.Range("I2").Value = "b"
$ws.Range("J2").Value = "html"

# Row 3-8: entity / variant names
$ws.Range("A3").Value = "uulm"
$ws.Range("B3").Value = "black"
$ws.Range("A4").Value = "uulm"
$ws.Range("B4").Value = "light blue"
$ws.Range("A5").Value = "ma-wi"
$ws.Range("B5").Value = "green"
$ws.Range("A6").Value = "ing-inf"
$ws.Range("B6").Value = "red"
$ws.Range("A7").Value = "nat"
$ws.Range("B7").Value = "orange"
$ws.Range("A8").Value = "med"
$ws.Range("B8").Value = "blue"

# Base RGB values
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "#000000"

$ws.Range("C4").Value = 137
$ws.Range("D4").Value = 162
$ws.Range("E4").Value = 179
$ws.Range("F4").Formula = '=_xlfn.CONCAT("#",DEC2HEX(C4),DEC2HEX(D4),DEC2HEX(E4))'

$ws.Range("C5").Value = 86
$ws.Range("D5").Value = 170
$ws.Range("E5").Value = 28
$ws.Range("F5:F8").Formula = '=_xlfn.CONCAT("#",DEC2HEX(C5),DEC2HEX(D5),DEC2HEX(E5))'

$ws.Range("C6").Value = 136
$ws.Range("D6").Value = 38
$ws.Range("E6").Value = 56

$ws.Range("C7").Value = 223
$ws.Range("D7").Value = 109
$ws.Range("E7").Value = 7
$ws.Range("F7").Formula = '=_xlfn.CONCAT("#",DEC2HEX(C7),DEC2HEX(D7),"0",DEC2HEX(E7))'

$ws.Range("C8").Value = 38
$ws.Range("D8").Value = 84
$ws.Range("E8").Value = 124

# Lightened RGB (mix with white by G1 percentage) + resulting html
$ws.Range("G3:G7").Formula = '=255*(1-$G$1)+C3*$G$1'
$ws.Range("H3:H7").Formula = '=255*(1-$G$1)+D3*$G$1'
$ws.Range("I3:I7").Formula = '=255*(1-$G$1)+E3*$G$1'
$ws.Range("G3:I7").NumberFormat = "0"

$ws.Range("G8").Formula = '=255*(1-$G$1)+C8*$G$1'
$ws.Range("H8:I8").Formula = '=255*(1-$G$1)+D8*$G$1'
$ws.Range("G8:I8").NumberFormat = "0"

$ws.Range("J3").Formula = '=_xlfn.CONCAT("#",DEC2HEX(G3),DEC2HEX(H3),DEC2HEX(I3))'
$ws.Range("J4:J8").Formula = '=_xlfn.CONCAT("#",DEC2HEX(G4),DEC2HEX(H4),DEC2HEX(I4))'

# Column widths (best-fit approximations)
$ws.Range("A1:A8").EntireColumn.AutoFit()
$ws.Range("B1:B8").EntireColumn.AutoFit()
$ws.Range("F1:F8").EntireColumn.AutoFit()
$ws.Range("J1:J8").EntireColumn.AutoFit()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Activate()
$ws.Range("G2").Select()

Write-Host "uulm sheet built"
